# "Add files via upload" edit for MeetTest_ScoreSheet.xlsx (Sheet1):
#   - The `effnet` rows (13, 16, 20) had an unknown learning rate recorded
#     as the placeholder text "?" in column C ("lr"). Fill in the real
#     value, 5E-04, formatted the same way as the rest of that column
#     (scientific notation, e.g. 9.00E-05).
#   - Move the saved active-cell selection from H11 to H15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the number format already applied to the other "lr" cells (column C)
# so the newly-filled cells render identically to their neighbours.
$lrNumberFormat = $ws.Range("C2").NumberFormat

foreach ($addr in @("C13", "C16", "C20")) {
    $cell = $ws.Range($addr)
    $cell.Value = 0.0005
    $cell.NumberFormat = $lrNumberFormat
}

# Update the sheet's active cell/selection (was H11).
[void]$ws.Range("H15").Select()
